# Controladora.xlsx - "Unidade de controle finalizada"
#
# 1. The ALUCtrl column (originally column F) is moved so it becomes the
#    last of the control-signal columns, after DataToReg: the old
#    F/G/H/I order (ALUCtrl, Branch, DataWrite, DataToReg) becomes
#    G/H/I/F -> new F/G/H/I order (Branch, DataWrite, DataToReg, ALUCtrl).
#    Doing this as a real column cut+insert also carries the cell
#    formatting (style) along with the data, and shifts only columns
#    F..I (columns further right, like the notes table in M/N, are left
#    untouched because the insert point is column J).
# 2. The active selection ends up on J28.
# 3. Page setup is switched to A4 paper, portrait orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Move column F (ALUCtrl) to after column I (DataToReg) ---------
$ws.Columns("F:F").Cut()
$ws.Columns("J:J").Insert()

# --- 2. Selection -------------------------------------------------------
$ws.Range("J28").Select()

# --- 3. Page setup: A4, portrait ----------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
